$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as TEXT (not auto-converted
# to a number) by building it as a text formula in a scratch cell, then pasting
# only the *value* back onto the target cell. This avoids Excel applying a
# quote-prefix/text number-format to the target cell (keeps its style untouched),
# matching how the source file stores these as plain (unstyled) text cells.
function Set-TextValue($range, $val) {
    $scratch = $ws.Range("Z1")
    $scratch.Formula = '="' + $val + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

# Row 2
$ws.Range("D2").Value = "29.018.83"
$ws.Range("E2").Value = "  -2.00%  "

# Row 3
$ws.Range("D3").Value = "1.964.86"
$ws.Range("E3").Value = "  -2.15%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.009"
$ws.Range("E4").Value = "  -0.35%  "

# Row 5
Set-TextValue $ws.Range("D5") "327.83"
$ws.Range("E5").Value = "  -0.90%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.007"
$ws.Range("E6").Value = "  -0.37%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4970"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.4204"
$ws.Range("E8").Value = "  -0.32%  "

# Row 9
Set-TextValue $ws.Range("D9") "52.93"
$ws.Range("E9").Value = "  -1.70%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.09207"
$ws.Range("E10").Value = "  +4.12%  "

# Row 11
Set-TextValue $ws.Range("D11") "1.097"
$ws.Range("E11").Value = "  -2.37%  "

# Row 12
Set-TextValue $ws.Range("D12") "22.79"
$ws.Range("E12").Value = "  -1.01%  "

# Row 13
$ws.Range("D13").Value = "1.962.38"
$ws.Range("E13").Value = "  +0.18%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.877"
$ws.Range("E14").Value = "  -3.17%  "

# Row 15
Set-TextValue $ws.Range("D15") "6.459"
$ws.Range("E15").Value = "  -0.84%  "

# Row 16
$ws.Range("E16").Value = "  -0.48%  "

# Row 17
Set-TextValue $ws.Range("D17") "91.74"
$ws.Range("E17").Value = "  -5.00%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.00001101"
$ws.Range("E18").Value = "  -0.62%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.06723"
$ws.Range("E19").Value = "  +1.52%  "

# Row 20
Set-TextValue $ws.Range("D20") "19.24"
$ws.Range("E20").Value = "  -1.69%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.007"
$ws.Range("E21").Value = "  -0.26%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.955"
$ws.Range("E22").Value = "  -0.90%  "

# Row 23
$ws.Range("D23").Value = "29.036.12"
$ws.Range("E23").Value = "  -1.99%  "

# Row 24
Set-TextValue $ws.Range("D24") "12.01"
$ws.Range("E24").Value = "  +0.73%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.267"
$ws.Range("E25").Value = "  -0.74%  "

# Row 26
$ws.Range("D26").Value = "2.223.19"
$ws.Range("E26").Value = "  +0.20%  "

# Row 27
Set-TextValue $ws.Range("D27") "20.63"
$ws.Range("E27").Value = "  -0.29%  "

# Row 28
Set-TextValue $ws.Range("D28") "155.73"
$ws.Range("E28").Value = "  -1.51%  "

# Row 29
Set-TextValue $ws.Range("D29") "6.380"
$ws.Range("E29").Value = "  -2.90%  "

# Row 30
$ws.Range("E30").Value = "  -3.86%  "

# Row 31
Set-TextValue $ws.Range("D31") "126.60"
$ws.Range("E31").Value = "  -0.47%  "

# Row 32
Set-TextValue $ws.Range("D32") "1.045"
$ws.Range("E32").Value = "  -1.33%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.09838"
$ws.Range("E33").Value = "  -1.38%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.519"
$ws.Range("E34").Value = "  -2.35%  "

# Row 35
Set-TextValue $ws.Range("D35") "5.823"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.714"
$ws.Range("E36").Value = "  -2.33%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.02432"
$ws.Range("E37").Value = "  -1.28%  "

# Row 38
$ws.Range("E38").Value = "  +1.82%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.06367"
$ws.Range("E39").Value = "  -0.56%  "

# Row 40
Set-TextValue $ws.Range("D40") "8.999"
$ws.Range("E40").Value = "  -6.47%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.6450"
$ws.Range("E41").Value = "  -1.11%  "

# Row 42
Set-TextValue $ws.Range("D42") "11.42"
$ws.Range("E42").Value = "  -3.49%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.1984"
$ws.Range("E43").Value = "  -4.44%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.007"
$ws.Range("E44").Value = "  -0.37%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.6208"
$ws.Range("E45").Value = "  -2.14%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.343"
$ws.Range("E46").Value = "  +5.77%  "

# Row 47
Set-TextValue $ws.Range("D47") "2.197"
$ws.Range("E47").Value = "  -1.50%  "

# Row 48
$ws.Range("E48").Value = "  -1.87%  "

# Row 49
Set-TextValue $ws.Range("D49") "3.471"
$ws.Range("E49").Value = "  -2.54%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.00000000325"
$ws.Range("E50").Value = "  -0.96%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.06971"
$ws.Range("E51").Value = "  -0.76%"

